$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new week's pair of records (Primera / Segunda) is
# inserted at the top of the data block (row 86), pushing all the
# existing rows from 86..169 down to 88..171 (dimension grows from
# A1:R169 to A1:R171).
$ws.Range("A86:A87").EntireRow.Insert()

# Populate the newly inserted row 86 (Primera).
$ws.Cells.Item(86,1).Value  = 11
$ws.Cells.Item(86,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(86,3).Value  = "Bíobío"
$ws.Cells.Item(86,4).Value  = 44827
$ws.Cells.Item(86,5).Value  = 8
$ws.Cells.Item(86,6).Value  = 100112044
$ws.Cells.Item(86,7).Value  = "Perejil"
$ws.Cells.Item(86,8).Value  = "Sin especificar"
$ws.Cells.Item(86,9).Value  = "Primera"
$ws.Cells.Item(86,10).Value = 200
$ws.Cells.Item(86,11).Value = 700
$ws.Cells.Item(86,12).Value = 800
$ws.Cells.Item(86,13).Value = 750
$ws.Cells.Item(86,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(86,15).Value = "Región Metropolitana"
$ws.Cells.Item(86,16).Value = 750
$ws.Cells.Item(86,17).Value = 1
$ws.Cells.Item(86,18).Value = "Hortaliza"

# Populate the newly inserted row 87 (Segunda).
$ws.Cells.Item(87,1).Value  = 11
$ws.Cells.Item(87,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(87,3).Value  = "Bíobío"
$ws.Cells.Item(87,4).Value  = 44827
$ws.Cells.Item(87,5).Value  = 8
$ws.Cells.Item(87,6).Value  = 100112044
$ws.Cells.Item(87,7).Value  = "Perejil"
$ws.Cells.Item(87,8).Value  = "Sin especificar"
$ws.Cells.Item(87,9).Value  = "Segunda"
$ws.Cells.Item(87,10).Value = 100
$ws.Cells.Item(87,11).Value = 600
$ws.Cells.Item(87,12).Value = 600
$ws.Cells.Item(87,13).Value = 600
$ws.Cells.Item(87,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(87,15).Value = "Región Metropolitana"
$ws.Cells.Item(87,16).Value = 600
$ws.Cells.Item(87,17).Value = 1
$ws.Cells.Item(87,18).Value = "Hortaliza"
